$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2, shifting existing data (rows 2-8) down to rows 3-9
$ws.Rows.Item(2).Insert()

# Populate the new row 2 with the new data (order matches shared-string insertion order)
$ws.Range("A2").Value = "Home"
$ws.Range("C2").Value = "Page"
$ws.Range("D2").Value = "//"
$ws.Range("B2").Value = "Now"

# Update selection to match the diff (D2 selected)
$ws.Range("D2").Select()
